$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 76, shifting existing rows 76-78 down to 77-79
$ws.Rows.Item(76).Insert()

# Populate the newly inserted row 76 with the new record
$ws.Cells.Item(76, 1).Value = 6
$ws.Cells.Item(76, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(76, 3).Value = "Metropolitana"
$ws.Cells.Item(76, 4).Value = 44610
$ws.Cells.Item(76, 5).Value = 13
$ws.Cells.Item(76, 6).Value = "Fruta"
$ws.Cells.Item(76, 7).Value = 100101
$ws.Cells.Item(76, 8).Value = "Berries"
$ws.Cells.Item(76, 9).Value = 100101008
$ws.Cells.Item(76, 10).Value = "Mora"
$ws.Cells.Item(76, 11).Value = "Sin especificar"
$ws.Cells.Item(76, 12).Value = "Primera"
$ws.Cells.Item(76, 13).Value = 250
$ws.Cells.Item(76, 14).Value = 6000
$ws.Cells.Item(76, 15).Value = 6000
$ws.Cells.Item(76, 16).Value = 6000
$ws.Cells.Item(76, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(76, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(76, 19).Value = 3000
$ws.Cells.Item(76, 20).Value = 2
